# Update scripts with new TPM-derived values for the Ptdss1-Jmjd6 ligand-receptor
# pair sheet (Ligand/Receptor expression values and derived specificity/weight
# metrics recomputed from the new TPM table).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 7).Value = 8.542726333333333
$ws.Cells.Item(2, 8).Value = 25.628179
$ws.Cells.Item(2, 9).Value = 0.2146499313812649
$ws.Cells.Item(2, 10).Value = 0.2146499313812649
$ws.Cells.Item(2, 13).Value = 10.64562733333334
$ws.Cells.Item(2, 14).Value = 31.936882
$ws.Cells.Item(2, 15).Value = 0.270105821029879
$ws.Cells.Item(2, 16).Value = 0.2701058210298791
$ws.Cells.Item(2, 17).Value = 90.94268095531979
$ws.Cells.Item(2, 18).Value = 818.484128597878
$ws.Cells.Item(2, 19).Value = 0.05797819594974375
$ws.Cells.Item(2, 20).Value = 0.05797819594974376
$ws.Cells.Item(3, 7).Value = 8.542726333333333
$ws.Cells.Item(3, 8).Value = 25.628179
$ws.Cells.Item(3, 9).Value = 0.2146499313812649
$ws.Cells.Item(3, 10).Value = 0.2146499313812649
$ws.Cells.Item(3, 15).Value = 0.2681760419860866
$ws.Cells.Item(3, 16).Value = 0.2681760419860866
$ws.Cells.Item(3, 17).Value = 90.2929382758591
$ws.Cells.Item(3, 18).Value = 812.6364444827319
$ws.Cells.Item(3, 19).Value = 0.0575639690104127
$ws.Cells.Item(3, 20).Value = 0.05756396901041271
$ws.Cells.Item(4, 7).Value = 8.542726333333333
$ws.Cells.Item(4, 8).Value = 25.628179
$ws.Cells.Item(4, 9).Value = 0.2146499313812649
$ws.Cells.Item(4, 10).Value = 0.2146499313812649
$ws.Cells.Item(4, 13).Value = 8.451328333333334
$ws.Cells.Item(4, 14).Value = 25.353985
$ws.Cells.Item(4, 15).Value = 0.2144310435440829
$ws.Cells.Item(4, 16).Value = 0.2144310435440829
$ws.Cells.Item(4, 17).Value = 72.19738510481278
$ws.Cells.Item(4, 18).Value = 649.776465943315
$ws.Cells.Item(4, 19).Value = 0.04602760878275042
$ws.Cells.Item(4, 20).Value = 0.04602760878275042
$ws.Cells.Item(5, 7).Value = 8.542726333333333
$ws.Cells.Item(5, 8).Value = 25.628179
$ws.Cells.Item(5, 9).Value = 0.2146499313812649
$ws.Cells.Item(5, 10).Value = 0.2146499313812649
$ws.Cells.Item(5, 13).Value = 4.945044
$ws.Cells.Item(5, 14).Value = 14.835132
$ws.Cells.Item(5, 15).Value = 0.1254679623686066
$ws.Cells.Item(5, 16).Value = 0.1254679623686067
$ws.Cells.Item(5, 17).Value = 42.244157598292
$ws.Cells.Item(5, 18).Value = 380.197418384628
$ws.Cells.Item(5, 19).Value = 0.02693168951296854
$ws.Cells.Item(5, 20).Value = 0.02693168951296855
$ws.Cells.Item(6, 7).Value = 8.542726333333333
$ws.Cells.Item(6, 8).Value = 25.628179
$ws.Cells.Item(6, 9).Value = 0.2146499313812649
$ws.Cells.Item(6, 10).Value = 0.2146499313812649
$ws.Cells.Item(6, 13).Value = 4.801233333333333
$ws.Cells.Item(6, 14).Value = 14.4037
$ws.Cells.Item(6, 15).Value = 0.1218191310713446
$ws.Cells.Item(6, 16).Value = 0.1218191310713447
$ws.Cells.Item(6, 17).Value = 41.01562242914444
$ws.Cells.Item(6, 18).Value = 369.1406018622999
$ws.Cells.Item(6, 19).Value = 0.02614846812538944
$ws.Cells.Item(6, 20).Value = 0.02614846812538945
$ws.Cells.Item(7, 9).Value = 0.2649602054889376
$ws.Cells.Item(7, 10).Value = 0.2649602054889376
$ws.Cells.Item(7, 13).Value = 10.64562733333334
$ws.Cells.Item(7, 14).Value = 31.936882
$ws.Cells.Item(7, 15).Value = 0.270105821029879
$ws.Cells.Item(7, 16).Value = 0.2701058210298791
$ws.Cells.Item(7, 17).Value = 112.2580905504058
$ws.Cells.Item(7, 18).Value = 1010.322814953652
$ws.Cells.Item(7, 19).Value = 0.07156729384383495
$ws.Cells.Item(7, 20).Value = 0.07156729384383496
$ws.Cells.Item(8, 9).Value = 0.2649602054889376
$ws.Cells.Item(8, 10).Value = 0.2649602054889376
$ws.Cells.Item(8, 15).Value = 0.2681760419860866
$ws.Cells.Item(8, 16).Value = 0.2681760419860866
$ws.Cells.Item(8, 19).Value = 0.07105597919184346
$ws.Cells.Item(8, 20).Value = 0.07105597919184348
$ws.Cells.Item(9, 9).Value = 0.2649602054889376
$ws.Cells.Item(9, 10).Value = 0.2649602054889376
$ws.Cells.Item(9, 13).Value = 8.451328333333334
$ws.Cells.Item(9, 14).Value = 25.353985
$ws.Cells.Item(9, 15).Value = 0.2144310435440829
$ws.Cells.Item(9, 16).Value = 0.2144310435440829
$ws.Cells.Item(9, 17).Value = 89.11921783546778
$ws.Cells.Item(9, 18).Value = 802.07296051921
$ws.Cells.Item(9, 19).Value = 0.05681569336064753
$ws.Cells.Item(9, 20).Value = 0.05681569336064755
$ws.Cells.Item(10, 9).Value = 0.2649602054889376
$ws.Cells.Item(10, 10).Value = 0.2649602054889376
$ws.Cells.Item(10, 13).Value = 4.945044
$ws.Cells.Item(10, 14).Value = 14.835132
$ws.Cells.Item(10, 15).Value = 0.1254679623686066
$ws.Cells.Item(10, 16).Value = 0.1254679623686067
$ws.Cells.Item(10, 17).Value = 52.145465903128
$ws.Cells.Item(10, 18).Value = 469.3091931281519
$ws.Cells.Item(10, 19).Value = 0.03324401709146431
$ws.Cells.Item(10, 20).Value = 0.03324401709146432
$ws.Cells.Item(11, 9).Value = 0.2649602054889376
$ws.Cells.Item(11, 10).Value = 0.2649602054889376
$ws.Cells.Item(11, 13).Value = 4.801233333333333
$ws.Cells.Item(11, 14).Value = 14.4037
$ws.Cells.Item(11, 15).Value = 0.1218191310713446
$ws.Cells.Item(11, 16).Value = 0.1218191310713447
$ws.Cells.Item(11, 17).Value = 50.62898309424444
$ws.Cells.Item(11, 18).Value = 455.6608478482
$ws.Cells.Item(11, 19).Value = 0.03227722200114731
$ws.Cells.Item(11, 20).Value = 0.03227722200114731
$ws.Cells.Item(12, 7).Value = 10.06002866666667
$ws.Cells.Item(12, 8).Value = 30.180086
$ws.Cells.Item(12, 9).Value = 0.2527746270611218
$ws.Cells.Item(12, 10).Value = 0.2527746270611218
$ws.Cells.Item(12, 13).Value = 10.64562733333334
$ws.Cells.Item(12, 14).Value = 31.936882
$ws.Cells.Item(12, 15).Value = 0.270105821029879
$ws.Cells.Item(12, 16).Value = 0.2701058210298791
$ws.Cells.Item(12, 17).Value = 107.0953161479836
$ws.Cells.Item(12, 18).Value = 963.8578453318521
$ws.Cells.Item(12, 19).Value = 0.06827589817786578
$ws.Cells.Item(12, 20).Value = 0.06827589817786579
$ws.Cells.Item(13, 7).Value = 10.06002866666667
$ws.Cells.Item(13, 8).Value = 30.180086
$ws.Cells.Item(13, 9).Value = 0.2527746270611218
$ws.Cells.Item(13, 10).Value = 0.2527746270611218
$ws.Cells.Item(13, 15).Value = 0.2681760419860866
$ws.Cells.Item(13, 16).Value = 0.2681760419860866
$ws.Cells.Item(13, 17).Value = 106.3301704876542
$ws.Cells.Item(13, 18).Value = 956.971534388888
$ws.Cells.Item(13, 19).Value = 0.06778809899976078
$ws.Cells.Item(13, 20).Value = 0.06778809899976079
$ws.Cells.Item(14, 7).Value = 10.06002866666667
$ws.Cells.Item(14, 8).Value = 30.180086
$ws.Cells.Item(14, 9).Value = 0.2527746270611218
$ws.Cells.Item(14, 10).Value = 0.2527746270611218
$ws.Cells.Item(14, 13).Value = 8.451328333333334
$ws.Cells.Item(14, 14).Value = 25.353985
$ws.Cells.Item(14, 15).Value = 0.2144310435440829
$ws.Cells.Item(14, 16).Value = 0.2144310435440829
$ws.Cells.Item(14, 17).Value = 85.02060530474556
$ws.Cells.Item(14, 18).Value = 765.18544774271
$ws.Cells.Item(14, 19).Value = 0.05420272706218272
$ws.Cells.Item(14, 20).Value = 0.05420272706218273
$ws.Cells.Item(15, 7).Value = 10.06002866666667
$ws.Cells.Item(15, 8).Value = 30.180086
$ws.Cells.Item(15, 9).Value = 0.2527746270611218
$ws.Cells.Item(15, 10).Value = 0.2527746270611218
$ws.Cells.Item(15, 13).Value = 4.945044
$ws.Cells.Item(15, 14).Value = 14.835132
$ws.Cells.Item(15, 15).Value = 0.1254679623686066
$ws.Cells.Item(15, 16).Value = 0.1254679623686067
$ws.Cells.Item(15, 17).Value = 49.747284397928
$ws.Cells.Item(15, 18).Value = 447.725559581352
$ws.Cells.Item(15, 19).Value = 0.03171511739584341
$ws.Cells.Item(15, 20).Value = 0.03171511739584341
$ws.Cells.Item(16, 7).Value = 10.06002866666667
$ws.Cells.Item(16, 8).Value = 30.180086
$ws.Cells.Item(16, 9).Value = 0.2527746270611218
$ws.Cells.Item(16, 10).Value = 0.2527746270611218
$ws.Cells.Item(16, 13).Value = 4.801233333333333
$ws.Cells.Item(16, 14).Value = 14.4037
$ws.Cells.Item(16, 15).Value = 0.1218191310713446
$ws.Cells.Item(16, 16).Value = 0.1218191310713447
$ws.Cells.Item(16, 17).Value = 48.30054496868888
$ws.Cells.Item(16, 18).Value = 434.7049047182
$ws.Cells.Item(16, 19).Value = 0.03079278542546906
$ws.Cells.Item(16, 20).Value = 0.03079278542546907
$ws.Cells.Item(17, 7).Value = 3.225032333333334
$ws.Cells.Item(17, 8).Value = 9.675097000000001
$ws.Cells.Item(17, 9).Value = 0.08103419705149875
$ws.Cells.Item(17, 10).Value = 0.08103419705149875
$ws.Cells.Item(17, 13).Value = 10.64562733333334
$ws.Cells.Item(17, 14).Value = 31.936882
$ws.Cells.Item(17, 15).Value = 0.270105821029879
$ws.Cells.Item(17, 16).Value = 0.2701058210298791
$ws.Cells.Item(17, 17).Value = 34.33249235861712
$ws.Cells.Item(17, 18).Value = 308.9924312275541
$ws.Cells.Item(17, 19).Value = 0.02188780832609207
$ws.Cells.Item(17, 20).Value = 0.02188780832609208
$ws.Cells.Item(18, 7).Value = 3.225032333333334
$ws.Cells.Item(18, 8).Value = 9.675097000000001
$ws.Cells.Item(18, 9).Value = 0.08103419705149875
$ws.Cells.Item(18, 10).Value = 0.08103419705149875
$ws.Cells.Item(18, 15).Value = 0.2681760419860866
$ws.Cells.Item(18, 16).Value = 0.2681760419860866
$ws.Cells.Item(18, 17).Value = 34.08720284940845
$ws.Cells.Item(18, 18).Value = 306.7848256446761
$ws.Cells.Item(18, 19).Value = 0.02173143023079154
$ws.Cells.Item(18, 20).Value = 0.02173143023079155
$ws.Cells.Item(19, 7).Value = 3.225032333333334
$ws.Cells.Item(19, 8).Value = 9.675097000000001
$ws.Cells.Item(19, 9).Value = 0.08103419705149875
$ws.Cells.Item(19, 10).Value = 0.08103419705149875
$ws.Cells.Item(19, 13).Value = 8.451328333333334
$ws.Cells.Item(19, 14).Value = 25.353985
$ws.Cells.Item(19, 15).Value = 0.2144310435440829
$ws.Cells.Item(19, 16).Value = 0.2144310435440829
$ws.Cells.Item(19, 17).Value = 27.25580713461612
$ws.Cells.Item(19, 18).Value = 245.302264211545
$ws.Cells.Item(19, 19).Value = 0.01737624743650972
$ws.Cells.Item(19, 20).Value = 0.01737624743650973
$ws.Cells.Item(20, 7).Value = 3.225032333333334
$ws.Cells.Item(20, 8).Value = 9.675097000000001
$ws.Cells.Item(20, 9).Value = 0.08103419705149875
$ws.Cells.Item(20, 10).Value = 0.08103419705149875
$ws.Cells.Item(20, 13).Value = 4.945044
$ws.Cells.Item(20, 14).Value = 14.835132
$ws.Cells.Item(20, 15).Value = 0.1254679623686066
$ws.Cells.Item(20, 16).Value = 0.1254679623686067
$ws.Cells.Item(20, 17).Value = 15.947926789756
$ws.Cells.Item(20, 18).Value = 143.531341107804
$ws.Cells.Item(20, 19).Value = 0.0101671955862277
$ws.Cells.Item(20, 20).Value = 0.0101671955862277
$ws.Cells.Item(21, 7).Value = 3.225032333333334
$ws.Cells.Item(21, 8).Value = 9.675097000000001
$ws.Cells.Item(21, 9).Value = 0.08103419705149875
$ws.Cells.Item(21, 10).Value = 0.08103419705149875
$ws.Cells.Item(21, 13).Value = 4.801233333333333
$ws.Cells.Item(21, 14).Value = 14.4037
$ws.Cells.Item(21, 15).Value = 0.1218191310713446
$ws.Cells.Item(21, 16).Value = 0.1218191310713447
$ws.Cells.Item(21, 17).Value = 15.48413273987778
$ws.Cells.Item(21, 18).Value = 139.3571946589
$ws.Cells.Item(21, 19).Value = 0.009871515471877696
$ws.Cells.Item(21, 20).Value = 0.009871515471877699
$ws.Cells.Item(22, 7).Value = 7.425629000000001
$ws.Cells.Item(22, 8).Value = 22.276887
$ws.Cells.Item(22, 9).Value = 0.1865810390171769
$ws.Cells.Item(22, 10).Value = 0.1865810390171769
$ws.Cells.Item(22, 13).Value = 10.64562733333334
$ws.Cells.Item(22, 14).Value = 31.936882
$ws.Cells.Item(22, 15).Value = 0.270105821029879
$ws.Cells.Item(22, 16).Value = 0.2701058210298791
$ws.Cells.Item(22, 17).Value = 79.0504790495927
$ws.Cells.Item(22, 18).Value = 711.4543114463341
$ws.Cells.Item(22, 19).Value = 0.05039662473234245
$ws.Cells.Item(22, 20).Value = 0.05039662473234246
$ws.Cells.Item(23, 7).Value = 7.425629000000001
$ws.Cells.Item(23, 8).Value = 22.276887
$ws.Cells.Item(23, 9).Value = 0.1865810390171769
$ws.Cells.Item(23, 10).Value = 0.1865810390171769
$ws.Cells.Item(23, 15).Value = 0.2681760419860866
$ws.Cells.Item(23, 16).Value = 0.2681760419860866
$ws.Cells.Item(23, 17).Value = 78.48570055911068
$ws.Cells.Item(23, 18).Value = 706.3713050319961
$ws.Cells.Item(23, 19).Value = 0.05003656455327808
$ws.Cells.Item(23, 20).Value = 0.0500365645532781
$ws.Cells.Item(24, 7).Value = 7.425629000000001
$ws.Cells.Item(24, 8).Value = 22.276887
$ws.Cells.Item(24, 9).Value = 0.1865810390171769
$ws.Cells.Item(24, 10).Value = 0.1865810390171769
$ws.Cells.Item(24, 13).Value = 8.451328333333334
$ws.Cells.Item(24, 14).Value = 25.353985
$ws.Cells.Item(24, 15).Value = 0.2144310435440829
$ws.Cells.Item(24, 16).Value = 0.2144310435440829
$ws.Cells.Item(24, 17).Value = 62.75642876052168
$ws.Cells.Item(24, 18).Value = 564.807858844695
$ws.Cells.Item(24, 19).Value = 0.04000876690199248
$ws.Cells.Item(24, 20).Value = 0.04000876690199249
$ws.Cells.Item(25, 7).Value = 7.425629000000001
$ws.Cells.Item(25, 8).Value = 22.276887
$ws.Cells.Item(25, 9).Value = 0.1865810390171769
$ws.Cells.Item(25, 10).Value = 0.1865810390171769
$ws.Cells.Item(25, 13).Value = 4.945044
$ws.Cells.Item(25, 14).Value = 14.835132
$ws.Cells.Item(25, 15).Value = 0.1254679623686066
$ws.Cells.Item(25, 16).Value = 0.1254679623686067
$ws.Cells.Item(25, 17).Value = 36.72006213267601
$ws.Cells.Item(25, 18).Value = 330.480559194084
$ws.Cells.Item(25, 19).Value = 0.02340994278210267
$ws.Cells.Item(25, 20).Value = 0.02340994278210268
$ws.Cells.Item(26, 7).Value = 7.425629000000001
$ws.Cells.Item(26, 8).Value = 22.276887
$ws.Cells.Item(26, 9).Value = 0.1865810390171769
$ws.Cells.Item(26, 10).Value = 0.1865810390171769
$ws.Cells.Item(26, 13).Value = 4.801233333333333
$ws.Cells.Item(26, 14).Value = 14.4037
$ws.Cells.Item(26, 15).Value = 0.1218191310713446
$ws.Cells.Item(26, 16).Value = 0.1218191310713447
$ws.Cells.Item(26, 17).Value = 35.65217747576667
$ws.Cells.Item(26, 18).Value = 320.8695972819
$ws.Cells.Item(26, 19).Value = 0.02272914004746114
$ws.Cells.Item(26, 20).Value = 0.02272914004746115
